# "removed machine_name_references and started using ui_lookup"
#
# The machine-readable label "TTL Name/ UPI" in A4 is replaced by the more
# human-friendly UI label "TTL Name or UPI". The column is then widened
# (re-fit) to accommodate the new, longer text, and column A is selected
# (as happens after an Autofit-style "select column, fit width" action).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the old machine-oriented string for the new UI lookup label.
$ws.Range("A4").Value = "TTL Name or UPI"

# Re-fit column A's width now that its longest label changed.
$ws.Columns.Item(1).ColumnWidth = 14.65

# Leave the whole of column A selected, as the end-user would after
# adjusting it.
$ws.Range("A1:A1048576").Select()
